$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")

# New journal entries (dates, week number, duration, type, description)
$rows = @(
    @{ Row = 14; Date = "03/17/2021"; Semaine = 4; Temps = "30min"; Type = "Documentation"; Desc = "Introduction du S.M.A.R.T" },
    @{ Row = 15; Date = "03/24/2021"; Semaine = 5; Temps = "45min"; Type = "Exercice"; Desc = "Fin de l'option d'écriture dans des fichiers externe au code" },
    @{ Row = 16; Date = "03/25/2021"; Semaine = 5; Temps = "1h30"; Type = "Exercice"; Desc = "Fin de l'option de lecture dans des fichiers externe au code" },
    @{ Row = 17; Date = "03/25/2021"; Semaine = 5; Temps = "45min"; Type = "Documentation"; Desc = "Théorie teste " },
    @{ Row = 18; Date = "03/29/2021"; Semaine = 6; Temps = "45min"; Type = "Documentation"; Desc = "Théorie sur la documentation" },
    @{ Row = 19; Date = "04/02/2021"; Semaine = 6; Temps = "3h30"; Type = "Exercice"; Desc = "Création de la documentation" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).NumberFormat = "d-mmm"
    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 2).Value = $r.Semaine
    $ws.Cells.Item($row, 3).Value = $r.Temps
    $ws.Cells.Item($row, 4).Value = $r.Type
    $ws.Cells.Item($row, 5).Value = $r.Desc
}

# Final row: only date + description (end of version 1.0)
$ws.Cells.Item(20, 1).NumberFormat = "d-mmm"
$ws.Cells.Item(20, 1).Value = "04/02/2021"
$ws.Cells.Item(20, 5).Value = "Fin de la version 1.0"

$ws.Range("E20").Select()
